$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transitions")
Write-Output $ws.Name
